$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Revert "attempt to add 2 predictors": remove the extra hyp labels that
# were added in column E for rows 5-9, and drop the stray H11 cell.
$ws.Range("E5:E9").Clear()
$ws.Range("H11").Clear()

# Row 9 had been stretched to height 15 to fit the (now removed) styled
# label in E9; auto-fit it back down to the sheet's default row height.
$ws.Rows.Item(9).AutoFit()

# Restore the prior scroll position / selection on the sheet
# (topLeftCell A7, active cell B16).
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 7
$ws.Range("B16").Select()
